# Add a "Save" column (H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header cell H1, styled like the other header cells (copy formatting from G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill H2:H12 with 0 (plain numeric, no special style)
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 8).Value = 0
}
